$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (worksheet rows 44..70), columns A..M.
# Column A repeats the running index (same value as column B); columns C..L hold
# per-category vulnerability counts; column M holds the "Ignore" count.
$data = @(
    @(42,42,0,0,0,0,0,0,0,0,0,0,3),
    @(43,43,0,0,0,0,0,0,0,1,0,0,11),
    @(44,44,0,0,0,0,0,0,0,1,0,0,9),
    @(45,45,0,0,0,0,0,0,0,0,0,1,11),
    @(46,46,1,0,0,0,0,0,0,0,0,2,6),
    @(47,47,0,0,0,0,0,0,0,0,0,0,14),
    @(48,48,0,0,0,0,0,0,0,0,0,1,6),
    @(49,49,0,0,0,1,0,0,0,0,0,0,6),
    @(50,50,0,0,1,0,0,0,0,0,0,0,28),
    @(51,51,0,0,1,0,0,0,0,0,0,0,28),
    @(52,52,0,0,0,0,0,0,0,0,0,0,4),
    @(53,53,0,0,0,1,1,0,0,0,0,0,5),
    @(54,54,1,0,0,0,0,0,0,0,0,0,3),
    @(55,55,0,0,0,0,0,0,0,0,0,13,57),
    @(56,56,1,0,0,9,1,0,0,11,0,17,86),
    @(57,57,0,0,0,0,0,0,0,1,0,0,3),
    @(58,58,0,0,0,0,0,0,0,1,0,0,6),
    @(59,59,0,0,0,0,0,0,0,0,0,0,7),
    @(60,60,0,0,0,0,0,0,0,0,0,0,5),
    @(61,61,0,0,0,0,0,0,0,0,0,0,5),
    @(62,62,0,0,0,0,0,0,0,0,0,0,4),
    @(63,63,0,0,0,0,0,0,0,0,0,0,6),
    @(64,64,0,0,0,0,0,0,0,0,0,0,6),
    @(65,65,0,0,0,0,0,0,0,0,0,0,5),
    @(66,66,0,0,0,0,0,0,0,0,0,0,3),
    @(67,67,1,0,0,0,0,0,0,0,0,0,6),
    @(68,68,1,0,0,0,0,0,0,1,0,1,3)
)

$startRow = 44
$endRow = $startRow + $data.Count - 1

# Column A uses the same (bold + bordered + centered) style as the rest of the
# index column above it; reuse that existing style by copying it down instead
# of rebuilding it by hand.
$ws.Range("A43").Copy($ws.Range("A" + $startRow + ":A" + $endRow))

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
